$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: D1 label text stays "Estimate(hr)" (unchanged, no write needed) ---

# --- Row 16 is new: give it a literal ID value (not the shared ROW() formula) ---
$ws.Range("A16").Value = 15

# --- Rewrite the backlog table body (rows 2-16) with the new, lower-level C++/Vulkan project plan ---
$ws.Range("B2").Value = "Set up C++ IDE"
$ws.Range("C2").Value = "M"
$ws.Range("D2").Value = 2
$ws.Range("E2").Value = "Set up IDE for C++ development"
$ws.Range("F2").Value = "Can complie hello world.cpp"
$ws.Range("B3").Value = "Add libraries to project"
$ws.Range("C3").Value = "M"
$ws.Range("D3").Value = 2
$ws.Range("E3").Value = "Intigrate Vulkan SDK and windowing library "
$ws.Range("F3").Value = "Library can be used in C++ code"
$ws.Range("B4").Value = "Create window hello world"
$ws.Range("C4").Value = "M"
$ws.Range("D4").Value = 2
$ws.Range("E4").Value = "Create program to create coloured window"
$ws.Range("F4").Value = "Window with clear colour "
$ws.Range("B5").Value = "Have triangle show up"
$ws.Range("C5").Value = "M"
$ws.Range("D5").Value = 2
$ws.Range("E5").Value = "alter program to have triangle show up in window"
$ws.Range("F5").Value = "Window now features a tri colour triangle"
$ws.Range("B6").Value = "Load mesh "
$ws.Range("C6").Value = "M"
$ws.Range("D6").Value = 2
$ws.Range("E6").Value = "Have mesh file loaded into the program"
$ws.Range("F6").Value = "Window now shows abutraty loaded mesh"
$ws.Range("B7").Value = "Load svg "
$ws.Range("C7").Value = "M"
$ws.Range("D7").Value = 2
$ws.Range("E7").Value = "Have SVG into program"
$ws.Range("F7").Value = "Abutrary SVG displays on screen as mesh"
$ws.Range("B8").Value = "SVG algorithm"
$ws.Range("C8").Value = "M"
$ws.Range("D8").Value = 2
$ws.Range("E8").Value = "Have SVG have shaders apply to display curves"
$ws.Range("F8").Value = "SVGs have curves instead of straight edges"
$ws.Range("B9").Value = "Create scene"
$ws.Range("C9").Value = "S"
$ws.Range("D9").Value = 2
$ws.Range("E9").Value = "Place plane and cube in scene with SVG"
$ws.Range("F9").Value = "Display renders floor, cube and SVG"
$ws.Range("B10").Value = "Lighting pass"
$ws.Range("C10").Value = "S"
$ws.Range("D10").Value = 2
$ws.Range("E10").Value = "Add lighting pass"
$ws.Range("F10").Value = "Cube, SVG and plane are lit"
$ws.Range("B11").Value = "Shadow map"
$ws.Range("C11").Value = "S"
$ws.Range("D11").Value = 2
$ws.Range("E11").Value = "Add Shadow pass"
$ws.Range("F11").Value = "Have Cube and SVG cast shadows"
$ws.Range("B12").Value = "Shadow strenght variance"
$ws.Range("C12").Value = "C"
$ws.Range("D12").Value = 2
$ws.Range("E12").Value = "Have SVG not render "
$ws.Range("F12").Value = "Have SVG only cast shadow"
$ws.Range("B13").Value = "Spline animation"
$ws.Range("C13").Value = "C"
$ws.Range("D13").Value = 2
$ws.Range("E13").Value = "Allow for lerping splines"
$ws.Range("F13").Value = "Have spline change shape"
$ws.Range("B14").Value = "Have shadow strenght animate"
$ws.Range("C14").Value = "C"
$ws.Range("D14").Value = 2
$ws.Range("E14").Value = "Have shadow strenght change"
$ws.Range("F14").Value = "Have shadow strenght change"
$ws.Range("B15").Value = "base animation on light rotation"
$ws.Range("C15").Value = "C"
$ws.Range("D15").Value = 2
$ws.Range("E15").Value = "base animations on Slerping on light rotation"
$ws.Range("F15").Value = "Light rotation changes spline shape and shadow strenght"
$ws.Range("B16").Value = "Have light rotation be controllable"
$ws.Range("C16").Value = "C"
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = "Add controller based on key inputs to change light direction"
$ws.Range("F16").Value = "Light roation can be changed with time"

# --- Sheet view: drop the custom zoom, go back to 100%, and move the selection ---
$excel.ActiveWindow.Zoom = 100
$ws.Range("A15").Select()
